# Add Yang Chenhui's (007杨晨辉) new weekly progress entry (row 3)
# to the "007杨晨辉" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("007杨晨辉")

# Make this the active/selected sheet (matches tabSelected/activeTab move).
$ws.Activate()

# Date range (reuses existing shared string "20240429-20240505").
$ws.Range("A3").Value = "20240429-20240505"

# Week number entered as text (leading apostrophe -> quote-prefix style,
# same pattern used for the "10" in row 2 above).
$ws.Range("B3").Value = "'12"

# Progress notes for the week.
$ws.Range("C3").Value = "`n1.图像配准任务`n（1）icp、NDT配准算法学习`n（2）配准相关算法推导，如PCA,KD-Tree，OCTree等"
$ws.Range("C3").WrapText = $true

# Difficulties column - re-uses the standard "暂无" placeholder text.
$ws.Range("D3").Value = "1.暂无；`n2. 暂无；`n3. 暂无；"
$ws.Range("D3").WrapText = $true

# Next plan column.
$ws.Range("E3").Value = "编码测试"
$ws.Range("E3").WrapText = $true

# Row height to match the taller wrapped content.
$ws.Rows.Item(3).RowHeight = 83.25

# Select the whole new row, as left by the editing session.
$ws.Rows.Item(3).Select()
